# GorestTestData.xlsx update: bump the test-user email suffixes from "4"
# to "10" (names in column A stay the same, only the Email column / B
# changes). Hyperlink addresses (mailto:) are left untouched - only the
# displayed text of the cell changes, matching the original edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("userdata")

$ws.Range("B2").Value = "abu10@gmail.com"
$ws.Range("B3").Value = "pallu10@gmail.com"
$ws.Range("B4").Value = "zarina10@gmail.com"
$ws.Range("B5").Value = "tahira10@gmail.com"

# Leave the selection where the author ended up after editing (C5).
$ws.Range("C5").Select()
